# Update the "Metadata" worksheet of the ValueSet workbook:
#  - refresh the generation Date
#  - insert a "Jurisdiction" row (with an empty value) right after "Contact",
#    pushing Description/Purpose/Copyright/Immutable down by one row
#
# (the "Include from Statut validatio" sheet only shifts shared-string
#  indices as a side effect - its cell values are untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Preserve the formatting of the last existing row (14: Immutable) onto the
# new last row (15) before we start shuffling values down, so the new row
# doesn't end up with default/no style.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Shift rows 11-14 (Description, Purpose, Copyright, Immutable) down to
# 12-15, working from the bottom up so we never clobber a value before it
# has been copied onward.
$ws.Range("A15").Value = $ws.Range("A14").Value2
$ws.Range("B15").Value = $ws.Range("B14").Value2

$ws.Range("A14").Value = $ws.Range("A13").Value2
$ws.Range("B14").Value = $ws.Range("B13").Value2

$ws.Range("A13").Value = $ws.Range("A12").Value2
$ws.Range("B13").Value = $ws.Range("B12").Value2

$ws.Range("A12").Value = $ws.Range("A11").Value2
$ws.Range("B12").Value = $ws.Range("B11").Value2

# Row 11 becomes the new "Jurisdiction" property with an empty value.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Refresh the generation Date value.
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"
